# Update cryptocurrency price/volume symbol data to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'238.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'21.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.386"
$ws.Range("D4").Style = "Normal"
$ws.Range("D6").Value = "'6.473"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'3.346"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.7967"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'1.028"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1391"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07334"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03148"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.02981"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09244"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.001662"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'3.263"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04771"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0005717"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.006245"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.005083"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Value = "'0.0001502"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'0.0004206"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'3.919"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'2.204"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Value = "'0.1054"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Value = "'0.04084"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006919"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.003505"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.1039"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.009165"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005449"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Value = "'0.03772"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
$ws.Range("D49").Value = "'0.00002103"
$ws.Range("D49").Style = "Normal"
